$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) At the bottom of the document: drop the bold "Play Dark Vortex
#    Free Slot Game | Yggdrasil Gaming" paragraph entirely, and swap
#    the text of the trailing italic paragraph for the new image
#    prompt (keeping its italic run formatting untouched). Do this
#    *before* touching the top of the document so the still-unique
#    "Read our review..." text can only match the real target
#    paragraph.
# ------------------------------------------------------------------

# Locate the bold banner paragraph (last occurrence of this exact
# text - the real title at the top keeps the Heading1 paragraph
# style, this is the plain-styled/run-bold duplicate near the end)
# and delete the whole paragraph, including its paragraph mark.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Play Dark Vortex Free Slot Game | Yggdrasil Gaming" -and $p.Style.NameLocal -ne "Heading 1") {
        $p.Range.Delete()
        break
    }
}

$oldDesc = "Read our review of Dark Vortex, a 5-reels and 243-3,125 paylines slot game packed with unique features and scary theme. Play Dark Vortex free today!"
$newDesc = "Create an enticing feature image for Dark Vortex with the following specifications: Style: Cartoony Subject: A happy Maya warrior with glasses should be the main focus of the image. The warrior should be wearing a detailed headpiece and extravagant clothing that make them stand out from the dark and ominous background. They should be smiling and holding a glowing Vortex symbol to add more excitement to the image. Background: The background should reflect the ominous and mysterious atmosphere of the game. It should feature a portal to another dimension, with eerie purple and green hues permeating the scene. The portal should be slightly open, revealing glimpses of the otherworldly realm beyond. Overall Mood: The image should be striking and eye-catching, capturing the attention of potential players. It should convey the thrill and otherworldly allure of the Dark Vortex game, tempting players to dive into the game and uncover its secrets."

$d.Content.Find.Execute($oldDesc, $true, $false, $false, $false, $false, $true, 1, $false, $newDesc, 2) | Out-Null

# ------------------------------------------------------------------
# 2) Insert a new "Meta description" paragraph right after the title
#    paragraph (Heading1, "Play Dark Vortex Free Slot Game | Yggdrasil
#    Gaming"). First grow an empty paragraph after the title (so the
#    title and every other paragraph are left untouched), then pour
#    the real run layout into that empty paragraph via InsertXML: a
#    leading empty run, a bold "Meta description" run, then a plain
#    run with the rest of the sentence - matching the rest of the
#    document's body paragraphs.
# ------------------------------------------------------------------
$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Dark Vortex, a 5-reels and 243-3,125 paylines slot game packed with unique features and scary theme. Play Dark Vortex free today!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range
$metaRange.Collapse(1)
$metaRange.InsertXML($metaXml)
